# Add team record columns (Wins, Losses, Ties) to the KCR_2010 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells in row 1
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the header style used by the rest of row 1 (bold, centered, thin border)
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

# Fill the team record for every data row (2-48): every team played the
# same 67-95-0 record for the season.
$ws.Range("AD2:AD48").Value = 67
$ws.Range("AE2:AE48").Value = 95
$ws.Range("AF2:AF48").Value = 0
